# Scheduled runner update: refresh market-derived profit figures in Excalibur_Profits workbook
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 562
$ws.Range("I5").Value = 562
$ws.Range("K5").Value = 562
$ws.Range("M5").Value = -447
$ws.Range("H11").Value = 7.25
$ws.Range("I11").Value = 7.25
$ws.Range("K11").Value = 7.25
$ws.Range("M11").Value = 132.75
$ws.Range("H53").Value = 530
$ws.Range("I53").Value = 412.5
$ws.Range("J53").Value = 1000
$ws.Range("K53").Value = 412.5
$ws.Range("L53").Value = 1000
$ws.Range("M53").Value = 224.5
$ws.Range("N53").Value = -2274
$ws.Range("H69").Value = 9772.77
$ws.Range("J69").Value = 10004.7
$ws.Range("L69").Value = 30014.1
$ws.Range("N69").Value = -31762.1
$ws.Range("H70").Value = 12037
$ws.Range("J70").Value = 22100.4
$ws.Range("L70").Value = 66301.20000000001
$ws.Range("N70").Value = -66841.20000000001
$ws.Range("H72").Value = 9772.77
$ws.Range("J72").Value = 10004.7
$ws.Range("L72").Value = 90042.3
$ws.Range("N72").Value = -98778.3
$ws.Range("H73").Value = 12037
$ws.Range("J73").Value = 22100.4
$ws.Range("L73").Value = 66301.20000000001
$ws.Range("N73").Value = -68173.20000000001
$ws.Range("H87").Value = 76038.55499999999
$ws.Range("I87").Value = 48333.332
$ws.Range("J87").Value = 89891.164
$ws.Range("K87").Value = 48333.332
$ws.Range("L87").Value = 89891.164
$ws.Range("M87").Value = -47085.332
$ws.Range("N87").Value = -92387.164
$ws.Range("H90").Value = 76038.55499999999
$ws.Range("I90").Value = 48333.332
$ws.Range("J90").Value = 89891.164
$ws.Range("K90").Value = 144999.996
$ws.Range("L90").Value = 269673.492
$ws.Range("M90").Value = -138759.996
$ws.Range("N90").Value = -282153.492
$ws.Range("H115").Value = 2237
$ws.Range("J115").Value = 3000
$ws.Range("L115").Value = 9000
$ws.Range("N115").Value = -12134
$ws.Range("H127").Value = 2598
$ws.Range("I127").Value = 2000
$ws.Range("J127").Value = 2797.3333
$ws.Range("K127").Value = 6000
$ws.Range("L127").Value = 8391.999899999999
$ws.Range("M127").Value = -1040
$ws.Range("N127").Value = -18311.9999
$ws.Range("H132").Value = 295768.4
$ws.Range("I132").Value = 310709.78
$ws.Range("K132").Value = 932129.3400000001
$ws.Range("M132").Value = -929599.3400000001
$ws.Range("H138").Value = 3028.7542
$ws.Range("I138").Value = 2037.6897
$ws.Range("J138").Value = 3926.9062
$ws.Range("K138").Value = 6113.0691
$ws.Range("L138").Value = 11780.7186
$ws.Range("M138").Value = -973.0690999999997
$ws.Range("N138").Value = -22060.7186

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17622.092
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 17622.092
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 17622.092
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -18196.092
$ws.Range("H38").Value = 2847.3333
$ws.Range("I38").Value = 2500
$ws.Range("J38").Value = 3021
$ws.Range("K38").Value = 2500
$ws.Range("L38").Value = 3021
$ws.Range("M38").Value = -2033
$ws.Range("N38").Value = -3955
$ws.Range("H45").Value = 4073.8
$ws.Range("I45").Value = 3508.2307
$ws.Range("J45").Value = 7750
$ws.Range("K45").Value = 3508.2307
$ws.Range("L45").Value = 7750
$ws.Range("M45").Value = -3131.2307
$ws.Range("N45").Value = -8504
$ws.Range("H122").Value = 2193.4375
$ws.Range("I122").Value = 2238.4285
$ws.Range("K122").Value = 6715.2855
$ws.Range("M122").Value = -4265.2855
$ws.Range("H133").Value = 69998.5
$ws.Range("J133").Value = 69998.5
$ws.Range("L133").Value = 69998.5
$ws.Range("N133").Value = -75058.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3551.1738
$ws.Range("I105").Value = 3718.4546
$ws.Range("K105").Value = 3718.4546
$ws.Range("M105").Value = -1971.4546

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4642.3335
$ws.Range("I99").Value = 2542
$ws.Range("J99").Value = 6142.5713
$ws.Range("K99").Value = 2542
$ws.Range("L99").Value = 6142.5713
$ws.Range("M99").Value = -1044
$ws.Range("N99").Value = -9138.5713
$ws.Range("H126").Value = 4642.3335
$ws.Range("I126").Value = 2542
$ws.Range("J126").Value = 6142.5713
$ws.Range("K126").Value = 7626
$ws.Range("L126").Value = 18427.7139
$ws.Range("M126").Value = -5156
$ws.Range("N126").Value = -23367.7139

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 200030.36
$ws.Range("I2").Value = 294148
$ws.Range("J2").Value = 30.375
$ws.Range("K2").Value = 1764888
$ws.Range("L2").Value = 182.25
$ws.Range("M2").Value = -1764775
$ws.Range("N2").Value = -408.25
$ws.Range("H93").Value = 5827
$ws.Range("J93").Value = 6009.7
$ws.Range("L93").Value = 18029.1
$ws.Range("N93").Value = -21773.1
$ws.Range("H107").Value = 644.7857
$ws.Range("J107").Value = 666.0909
$ws.Range("L107").Value = 1998.2727
$ws.Range("N107").Value = -5838.2727
$ws.Range("H113").Value = 2960.3333
$ws.Range("J113").Value = 3000.3572
$ws.Range("L113").Value = 9001.071599999999
$ws.Range("N113").Value = -13341.0716
$ws.Range("H129").Value = 2538.9333
$ws.Range("J129").Value = 4265.143
$ws.Range("L129").Value = 12795.429
$ws.Range("N129").Value = -22795.429

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 150000
$ws.Range("J109").Value = 150000
$ws.Range("L109").Value = 150000
$ws.Range("N109").Value = -152080
$ws.Range("H132").Value = 3584.2632
$ws.Range("I132").Value = 2852.923
$ws.Range("K132").Value = 8558.769
$ws.Range("M132").Value = -6028.769
$ws.Range("H133").Value = 89892.336
$ws.Range("J133").Value = 89892.336
$ws.Range("L133").Value = 89892.336
$ws.Range("N133").Value = -100012.336

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4439.222
$ws.Range("I40").Value = 6281.353
$ws.Range("K40").Value = 6281.353
$ws.Range("M40").Value = -6145.353
$ws.Range("H88").Value = 10001
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 10001
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 10001
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -10857
$ws.Range("H91").Value = 10001
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 10001
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 10001
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -12965
$ws.Range("H122").Value = 3927.7827
$ws.Range("I122").Value = 3804.3333
$ws.Range("K122").Value = 11412.9999
$ws.Range("M122").Value = -8962.999899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3666.6667
$ws.Range("I81").Value = 3666.6667
$ws.Range("K81").Value = 7333.3334
$ws.Range("M81").Value = -6272.3334
$ws.Range("H84").Value = 3666.6667
$ws.Range("I84").Value = 3666.6667
$ws.Range("K84").Value = 36666.667
$ws.Range("M84").Value = -31362.667
$ws.Range("H132").Value = 11160966
$ws.Range("I132").Value = 961290.1
$ws.Range("K132").Value = 2883870.3
$ws.Range("M132").Value = -2881340.3
